$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10..114 shift down to 11..115.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44545
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112024
$ws.Cells.Item(10, 7).Value = "Choclo"
$ws.Cells.Item(10, 8).Value = "Choclero"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 10000
$ws.Cells.Item(10, 11).Value = 350
$ws.Cells.Item(10, 12).Value = 400
$ws.Cells.Item(10, 13).Value = 375
$ws.Cells.Item(10, 14).Value = "`$/unidad"
$ws.Cells.Item(10, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 16).Value = 375
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"
